$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer "Update automatically" date field: 18/10/2016 -> 19/10/2016
#    This cached field text lives on the slide master and on every
#    slide layout's Date placeholder (ppPlaceholderDate = 16).
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "18/10/2016") {
                $tr.Text = "19/10/2016"
            }
        }
    }
}

# Slide master footer date placeholder
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's footer date placeholder
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 1 "Su-rvival" -> "SU-rvival" (capitalise the U)
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(3)
$titleRange = $titleShape.TextFrame.TextRange
if ($titleRange.Text -eq "Su-rvival") {
    $titleRange.Characters(1, 3).Text = "SU-"
}
